$d = $word.ActiveDocument

# --- Change 1: expand mention of osgEarth-2.6 directory to mention osgEarth-2.7 as well ---
$old1 = "located in the osgEarth-2.6 directory. "
$new1 = "located in the osgEarth-2.6 directory for osgEarth version 2.6. An osgEarth-2.7 directory has been added which contains the updates for the feature driver required for osgEarth version 2.7."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Change 2: rewrite the OpenSceneGraph / CDB texture paragraph section ---
$old2 = "The OpenSceneGraph directory contains an update to the OpenFlight plugin for osg that is needed to load textures for CDB geospecific models that use the CDB 3.0 specification. In this case the model textures are found in a .zip file and require re-pathing  of the texture names from the highest "
$new2 = "The OpenSceneGraph-3.2.1 directory contains an update to the OpenFlight plugin for osg that is needed to load textures for CDB geospecific models that use the CDB 3.0 specification for OpenSceneGraph version 3.2.1. The OpenSceneGraph-3.4.0 directory has been added and contains the same changes merged into the OpenSceneGraph version 3.4.0 OpenFlight plugin. For the CDB  case the model textures are found in a .zip file and require re-pathing  of the texture names from the highest "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Remove the _GoBack bookmark that sat after "Update 28-May-2015:" ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Append the new "Update 19-September-2015" content at the end of the document ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $lastPara.Range

# blank paragraph
$r.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $lastPara.Range

# "Update 19-September-2015" paragraph
$r.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertAfter("Update 19-September-2015")
$r = $lastPara.Range

# Big update paragraph
$r.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$updateText = "Added updates for utilization with osgEarth 2.7 and OpenSceneGraph version 3.4.0. Updates for CMake and missing round function in Visual Studio versions prior to VS2013 were provided by Rapheal Cuisinier. Updates were added to allow for utilization of GDAL2.0 as well as maintain compatibility with previous versions of GDAL."
$lastPara.Range.InsertAfter($updateText)
$r = $lastPara.Range

# Final empty paragraph that holds the _GoBack bookmark
$r.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$nr = $lastPara.Range
$nr.Collapse(1)
$d.Bookmarks.Add("_GoBack", $nr) | Out-Null
